$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty "Unnamed" columns O:Z (12 columns, headers only, no
# data beneath them). This shifts the former "Imágenes" column (AA)
# left so it becomes column O.
$ws.Range("O1:Z1").EntireColumn.Delete()

# Fix a typo in the "Cárdigan Montanhas" image-URL list: the combining
# acute accent over the "a" was previously mis-encoded.
$cell = $ws.Range("O16")
$old = $cell.Value2
$new = $old -replace [regex]::Escape("C%C2%A0rdigan"), "Ca%CC%81rdigan"
$cell.Value2 = $new
